$wb = $excel.ActiveWorkbook

# --- Step 1: Major sheet - shift rows 19..40 up into 18..39 (B,C,D cols), drop old row 40 ---
$major = $wb.Worksheets.Item("Major")
$major.Cells.Item(18,2).Value = "Tadawul All Share"
$major.Cells.Item(18,3).Value = "Saudi Arabia"
$major.Cells.Item(18,4).Value = "SARUSD=X"
$major.Cells.Item(19,2).Value = "Nikkei 225"
$major.Cells.Item(19,3).Value = "Japan"
$major.Cells.Item(19,4).Value = "JPYUSD=X"
$major.Cells.Item(20,2).Value = "S&P/ASX 200"
$major.Cells.Item(20,3).Value = "Australia"
$major.Cells.Item(20,4).Value = "AUDUSD=X"
$major.Cells.Item(21,2).Value = "NZX 50"
$major.Cells.Item(21,3).Value = "New Zealand"
$major.Cells.Item(21,4).Value = "NZDUSD=X"
$major.Cells.Item(22,2).Value = "SZSE Component"
$major.Cells.Item(22,3).Value = "China"
$major.Cells.Item(22,4).Value = "CNYUSD=X"
$major.Cells.Item(23,2).Value = "China A50"
$major.Cells.Item(23,3).Value = "China"
$major.Cells.Item(23,4).Value = "CNYUSD=X"
$major.Cells.Item(24,2).Value = "Hang Seng"
$major.Cells.Item(24,3).Value = "Hong Kong"
$major.Cells.Item(24,4).Value = "HKDUSD=X"
$major.Cells.Item(25,2).Value = "Taiwan Weighted"
$major.Cells.Item(25,3).Value = "Taiwan"
$major.Cells.Item(25,4).Value = "TWDUSD=X"
$major.Cells.Item(26,2).Value = "SET"
$major.Cells.Item(26,3).Value = "Thailand"
$major.Cells.Item(26,4).Value = "THBUSD=X"
$major.Cells.Item(27,2).Value = "KOSPI"
$major.Cells.Item(27,3).Value = "South Korea"
$major.Cells.Item(27,4).Value = "KRWUSD=X"
$major.Cells.Item(28,2).Value = "IDX Composite"
$major.Cells.Item(28,3).Value = "Indonesia"
$major.Cells.Item(28,4).Value = "IDRUSD=X"
$major.Cells.Item(29,2).Value = "Nifty 50"
$major.Cells.Item(29,3).Value = "India"
$major.Cells.Item(29,4).Value = "INRUSD=X"
$major.Cells.Item(30,2).Value = "BSE Sensex"
$major.Cells.Item(30,3).Value = "India"
$major.Cells.Item(30,4).Value = "INRUSD=X"
$major.Cells.Item(31,2).Value = "PSEi Composite"
$major.Cells.Item(31,3).Value = "Philippines"
$major.Cells.Item(31,4).Value = "PHPUSD=X"
$major.Cells.Item(32,2).Value = "HNX 30"
$major.Cells.Item(32,3).Value = "Vietnam"
$major.Cells.Item(32,4).Value = "VNDUSD=X"
$major.Cells.Item(33,2).Value = "S&P CLX IPSA"
$major.Cells.Item(33,3).Value = "Chile"
$major.Cells.Item(33,4).Value = "CLPUSD=X"
$major.Cells.Item(34,2).Value = "COLCAP"
$major.Cells.Item(34,3).Value = "Colombia"
$major.Cells.Item(34,4).Value = "COPUSD=X"
$major.Cells.Item(35,2).Value = "KLCI"
$major.Cells.Item(35,3).Value = "Malaysia"
$major.Cells.Item(35,4).Value = "MYRUSD=X"
$major.Cells.Item(36,2).Value = "FTSE Singapore"
$major.Cells.Item(36,3).Value = "Singapore"
$major.Cells.Item(36,4).Value = "SGDUSD=X"
$major.Cells.Item(37,2).Value = "FTSE/JSE Top 40"
$major.Cells.Item(37,3).Value = "South Africa"
$major.Cells.Item(37,4).Value = "ZARUSD=X"
$major.Cells.Item(38,2).Value = "KOSDAQ"
$major.Cells.Item(38,3).Value = "South Korea"
$major.Cells.Item(38,4).Value = "KRWUSD=X"
$major.Cells.Item(39,2).Value = "S&P Lima General"
$major.Cells.Item(39,3).Value = "Peru"
$major.Cells.Item(39,4).Value = "PENUSD=X"
$major.Rows("40:40").Delete()
$major.Range("B1:C39").Select()

# --- Step 2: add "Sheet1" after Commodities, holding the displaced Turkey row ---
$afterCommodities = $wb.Worksheets.Item("Commodities")
$sheet1 = $wb.Worksheets.Add($null, $afterCommodities)
$sheet1.Name = "Sheet1"
$sheet1.Cells.Item(1,1).Formula = "=Major!A39+1"
$sheet1.Cells.Item(1,2).Value = "BIST 100"
$sheet1.Cells.Item(1,3).Value = "Turkey"
$sheet1.Cells.Item(1,4).Value = "TRYUSD=X"
$sheet1.Cells.Item(1,1).Borders.LineStyle = 1
$sheet1.Range("B1:D1").Borders.LineStyle = 1
$sheet1.Range("A1:C246").Select()

# --- Step 3: add "iso" sheet after Sheet1, with Lat/Lon/iso_alpha/iso_num lookup table ---
$iso = $wb.Worksheets.Add($null, $sheet1)
$iso.Name = "iso"
$iso.Cells.Item(1,1).Value = "Indices"
$iso.Cells.Item(1,2).Value = "Country"
$iso.Cells.Item(1,3).Value = "Lat"
$iso.Cells.Item(1,4).Value = "Lon"
$iso.Cells.Item(1,5).Value = "iso_alpha"
$iso.Cells.Item(1,6).Value = "iso_num"
$iso.Cells.Item(2,1).Value = "S&P 500"
$iso.Cells.Item(2,2).Value = "United States"
$iso.Cells.Item(2,3).Value = 37.090240000000001
$iso.Cells.Item(2,4).Value = -95.712890999999999
$iso.Cells.Item(2,5).Value = "USA"
$iso.Cells.Item(2,6).Value = 840
$iso.Cells.Item(3,1).Value = "Nasdaq"
$iso.Cells.Item(3,2).Value = "United States"
$iso.Cells.Item(3,3).Value = 37.090240000000001
$iso.Cells.Item(3,4).Value = -95.712890999999999
$iso.Cells.Item(3,5).Value = "USA"
$iso.Cells.Item(3,6).Value = 840
$iso.Cells.Item(4,1).Value = "SmallCap 2000"
$iso.Cells.Item(4,2).Value = "United States"
$iso.Cells.Item(4,3).Value = 37.090240000000001
$iso.Cells.Item(4,4).Value = -95.712890999999999
$iso.Cells.Item(4,5).Value = "USA"
$iso.Cells.Item(4,6).Value = 840
$iso.Cells.Item(5,1).Value = "S&P/TSX"
$iso.Cells.Item(5,2).Value = "Canada"
$iso.Cells.Item(5,3).Value = 56.130366000000002
$iso.Cells.Item(5,4).Value = -106.346771
$iso.Cells.Item(5,5).Value = "CAN"
$iso.Cells.Item(5,6).Value = 124
$iso.Cells.Item(6,1).Value = "Bovespa"
$iso.Cells.Item(6,2).Value = "Brazil"
$iso.Cells.Item(6,3).Value = -14.235004
$iso.Cells.Item(6,4).Value = -51.925280000000001
$iso.Cells.Item(6,5).Value = "BRA"
$iso.Cells.Item(6,6).Value = 76
$iso.Cells.Item(7,1).Value = "S&P/BMV IPC"
$iso.Cells.Item(7,2).Value = "Mexico"
$iso.Cells.Item(7,3).Value = 23.634501
$iso.Cells.Item(7,4).Value = -102.552784
$iso.Cells.Item(7,5).Value = "MEX"
$iso.Cells.Item(7,6).Value = 484
$iso.Cells.Item(8,1).Value = "DAX"
$iso.Cells.Item(8,2).Value = "Germany"
$iso.Cells.Item(8,3).Value = 51.165691000000002
$iso.Cells.Item(8,4).Value = 10.451525999999999
$iso.Cells.Item(8,5).Value = "DEU"
$iso.Cells.Item(8,6).Value = 276
$iso.Cells.Item(9,1).Value = "FTSE 100"
$iso.Cells.Item(9,2).Value = "United Kingdom"
$iso.Cells.Item(9,3).Value = 55.378050999999999
$iso.Cells.Item(9,4).Value = -3.4359730000000002
$iso.Cells.Item(9,5).Value = "GBR"
$iso.Cells.Item(9,6).Value = 826
$iso.Cells.Item(10,1).Value = "CAC 40"
$iso.Cells.Item(10,2).Value = "France"
$iso.Cells.Item(10,3).Value = 46.227637999999999
$iso.Cells.Item(10,4).Value = 2.213749
$iso.Cells.Item(10,5).Value = "FRA"
$iso.Cells.Item(10,6).Value = 250
$iso.Cells.Item(11,1).Value = "Euro Stoxx 50"
$iso.Cells.Item(11,2).Value = "Germany"
$iso.Cells.Item(11,3).Value = 51.165691000000002
$iso.Cells.Item(11,4).Value = 10.451525999999999
$iso.Cells.Item(11,5).Value = "DEU"
$iso.Cells.Item(11,6).Value = 276
$iso.Cells.Item(12,1).Value = "AEX"
$iso.Cells.Item(12,2).Value = "Netherlands"
$iso.Cells.Item(12,3).Value = 52.132632999999998
$iso.Cells.Item(12,4).Value = 5.2912660000000002
$iso.Cells.Item(12,5).Value = "NLD"
$iso.Cells.Item(12,6).Value = 528
$iso.Cells.Item(13,1).Value = "IBEX 35"
$iso.Cells.Item(13,2).Value = "Spain"
$iso.Cells.Item(13,3).Value = 40.463667000000001
$iso.Cells.Item(13,4).Value = -3.7492200000000002
$iso.Cells.Item(13,5).Value = "ESP"
$iso.Cells.Item(13,6).Value = 724
$iso.Cells.Item(14,1).Value = "FTSE MIB"
$iso.Cells.Item(14,2).Value = "Italy"
$iso.Cells.Item(14,3).Value = 41.871940000000002
$iso.Cells.Item(14,4).Value = 12.56738
$iso.Cells.Item(14,5).Value = "ITA"
$iso.Cells.Item(14,6).Value = 380
$iso.Cells.Item(15,1).Value = "SMI"
$iso.Cells.Item(15,2).Value = "Switzerland"
$iso.Cells.Item(15,3).Value = 46.818187999999999
$iso.Cells.Item(15,4).Value = 8.2275120000000008
$iso.Cells.Item(15,5).Value = "CHE"
$iso.Cells.Item(15,6).Value = 756
$iso.Cells.Item(16,1).Value = "OMXC25"
$iso.Cells.Item(16,2).Value = "Denmark"
$iso.Cells.Item(16,3).Value = 56.263919999999999
$iso.Cells.Item(16,4).Value = 9.5017849999999999
$iso.Cells.Item(16,5).Value = "DNK"
$iso.Cells.Item(16,6).Value = 208
$iso.Cells.Item(17,1).Value = "MOEX"
$iso.Cells.Item(17,2).Value = "Russia"
$iso.Cells.Item(17,3).Value = 61.524009999999997
$iso.Cells.Item(17,4).Value = 105.31875599999999
$iso.Cells.Item(17,5).Value = "RUS"
$iso.Cells.Item(17,6).Value = 643
$iso.Cells.Item(18,1).Value = "Tadawul All Share"
$iso.Cells.Item(18,2).Value = "Saudi Arabia"
$iso.Cells.Item(18,3).Value = 23.885942
$iso.Cells.Item(18,4).Value = 45.079161999999997
$iso.Cells.Item(18,5).Value = "SAU"
$iso.Cells.Item(18,6).Value = 682
$iso.Cells.Item(19,1).Value = "Nikkei 225"
$iso.Cells.Item(19,2).Value = "Japan"
$iso.Cells.Item(19,3).Value = 36.204824000000002
$iso.Cells.Item(19,4).Value = 138.25292400000001
$iso.Cells.Item(19,5).Value = "JPN"
$iso.Cells.Item(19,6).Value = 392
$iso.Cells.Item(20,1).Value = "S&P/ASX 200"
$iso.Cells.Item(20,2).Value = "Australia"
$iso.Cells.Item(20,3).Value = -25.274398000000001
$iso.Cells.Item(20,4).Value = 133.775136
$iso.Cells.Item(20,5).Value = "AUS"
$iso.Cells.Item(20,6).Value = 36
$iso.Cells.Item(21,1).Value = "NZX 50"
$iso.Cells.Item(21,2).Value = "New Zealand"
$iso.Cells.Item(21,3).Value = -40.900556999999999
$iso.Cells.Item(21,4).Value = 174.88597100000001
$iso.Cells.Item(21,5).Value = "NZL"
$iso.Cells.Item(21,6).Value = 554
$iso.Cells.Item(22,1).Value = "SZSE Component"
$iso.Cells.Item(22,2).Value = "China"
$iso.Cells.Item(22,3).Value = 35.861660000000001
$iso.Cells.Item(22,4).Value = 104.195397
$iso.Cells.Item(22,5).Value = "CHN"
$iso.Cells.Item(22,6).Value = 156
$iso.Cells.Item(23,1).Value = "China A50"
$iso.Cells.Item(23,2).Value = "China"
$iso.Cells.Item(23,3).Value = 35.861660000000001
$iso.Cells.Item(23,4).Value = 104.195397
$iso.Cells.Item(23,5).Value = "CHN"
$iso.Cells.Item(23,6).Value = 156
$iso.Cells.Item(24,1).Value = "Hang Seng"
$iso.Cells.Item(24,2).Value = "Hong Kong"
$iso.Cells.Item(24,3).Value = 22.396428
$iso.Cells.Item(24,4).Value = 114.109497
$iso.Cells.Item(24,5).Value = "HKG"
$iso.Cells.Item(24,6).Value = 344
$iso.Cells.Item(25,1).Value = "Taiwan Weighted"
$iso.Cells.Item(25,2).Value = "Taiwan"
$iso.Cells.Item(25,3).Value = 23.69781
$iso.Cells.Item(25,4).Value = 120.960515
$iso.Cells.Item(25,5).Value = "TWN"
$iso.Cells.Item(25,6).Value = 158
$iso.Cells.Item(26,1).Value = "SET"
$iso.Cells.Item(26,2).Value = "Thailand"
$iso.Cells.Item(26,3).Value = 15.870032
$iso.Cells.Item(26,4).Value = 100.992541
$iso.Cells.Item(26,5).Value = "THA"
$iso.Cells.Item(26,6).Value = 764
$iso.Cells.Item(27,1).Value = "KOSPI"
$iso.Cells.Item(27,2).Value = "South Korea"
$iso.Cells.Item(27,3).Value = 35.907756999999997
$iso.Cells.Item(27,4).Value = 127.76692199999999
$iso.Cells.Item(27,5).Value = "KOR"
$iso.Cells.Item(27,6).Value = 410
$iso.Cells.Item(28,1).Value = "IDX Composite"
$iso.Cells.Item(28,2).Value = "Indonesia"
$iso.Cells.Item(28,3).Value = -0.78927499999999995
$iso.Cells.Item(28,4).Value = 113.92132700000001
$iso.Cells.Item(28,5).Value = "IDN"
$iso.Cells.Item(28,6).Value = 360
$iso.Cells.Item(29,1).Value = "Nifty 50"
$iso.Cells.Item(29,2).Value = "India"
$iso.Cells.Item(29,3).Value = 20.593684
$iso.Cells.Item(29,4).Value = 78.962879999999998
$iso.Cells.Item(29,5).Value = "IND"
$iso.Cells.Item(29,6).Value = 356
$iso.Cells.Item(30,1).Value = "BSE Sensex"
$iso.Cells.Item(30,2).Value = "India"
$iso.Cells.Item(30,3).Value = 20.593684
$iso.Cells.Item(30,4).Value = 78.962879999999998
$iso.Cells.Item(30,5).Value = "IND"
$iso.Cells.Item(30,6).Value = 356
$iso.Cells.Item(31,1).Value = "PSEi Composite"
$iso.Cells.Item(31,2).Value = "Philippines"
$iso.Cells.Item(31,3).Value = 12.879721
$iso.Cells.Item(31,4).Value = 121.774017
$iso.Cells.Item(31,5).Value = "PHL"
$iso.Cells.Item(31,6).Value = 608
$iso.Cells.Item(32,1).Value = "HNX 30"
$iso.Cells.Item(32,2).Value = "Vietnam"
$iso.Cells.Item(32,3).Value = 14.058324000000001
$iso.Cells.Item(32,4).Value = 108.277199
$iso.Cells.Item(32,5).Value = "VNM"
$iso.Cells.Item(32,6).Value = 704
$iso.Cells.Item(33,1).Value = "S&P CLX IPSA"
$iso.Cells.Item(33,2).Value = "Chile"
$iso.Cells.Item(33,3).Value = -35.675147000000003
$iso.Cells.Item(33,4).Value = -71.542968999999999
$iso.Cells.Item(33,5).Value = "CHL"
$iso.Cells.Item(33,6).Value = 152
$iso.Cells.Item(34,1).Value = "COLCAP"
$iso.Cells.Item(34,2).Value = "Colombia"
$iso.Cells.Item(34,3).Value = 4.5708679999999999
$iso.Cells.Item(34,4).Value = -74.297332999999995
$iso.Cells.Item(34,5).Value = "COL"
$iso.Cells.Item(34,6).Value = 170
$iso.Cells.Item(35,1).Value = "KLCI"
$iso.Cells.Item(35,2).Value = "Malaysia"
$iso.Cells.Item(35,3).Value = 4.2104840000000001
$iso.Cells.Item(35,4).Value = 101.97576599999999
$iso.Cells.Item(35,5).Value = "MYS"
$iso.Cells.Item(35,6).Value = 458
$iso.Cells.Item(36,1).Value = "FTSE Singapore"
$iso.Cells.Item(36,2).Value = "Singapore"
$iso.Cells.Item(36,3).Value = 1.3520829999999999
$iso.Cells.Item(36,4).Value = 103.819836
$iso.Cells.Item(36,5).Value = "SGP"
$iso.Cells.Item(36,6).Value = 702
$iso.Cells.Item(37,1).Value = "FTSE/JSE Top 40"
$iso.Cells.Item(37,2).Value = "South Africa"
$iso.Cells.Item(37,3).Value = -30.559481999999999
$iso.Cells.Item(37,4).Value = 22.937505999999999
$iso.Cells.Item(37,5).Value = "ZAF"
$iso.Cells.Item(37,6).Value = 710
$iso.Cells.Item(38,1).Value = "KOSDAQ"
$iso.Cells.Item(38,2).Value = "South Korea"
$iso.Cells.Item(38,3).Value = 35.907756999999997
$iso.Cells.Item(38,4).Value = 127.76692199999999
$iso.Cells.Item(38,5).Value = "KOR"
$iso.Cells.Item(38,6).Value = 410
$iso.Cells.Item(39,1).Value = "S&P Lima General"
$iso.Cells.Item(39,2).Value = "Peru"
$iso.Cells.Item(39,3).Value = -9.1899669999999993
$iso.Cells.Item(39,4).Value = -75.015152
$iso.Cells.Item(39,5).Value = "PER"
$iso.Cells.Item(39,6).Value = 604
$iso.Range("A1:F39").Borders.LineStyle = 1
$iso.Columns.Item(1).ColumnWidth = 16.85546875
$iso.Columns.Item(2).ColumnWidth = 15.42578125
$iso.Range("I19").Select()

# --- Step 4: make iso the active/selected sheet (workbook activeTab -> 5) ---
$iso.Activate()
